$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that get numeric-looking text stay as literal
# text (matching the source inlineStr cells) instead of being auto-coerced
# to numbers (which would silently drop things like trailing zeros).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.841.59'
$ws.Range('E2').Value = '  +3.28%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.673.60'
$ws.Range('E3').Value = '  +2.81%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '202.26'
$ws.Range('E5').Value = '  +10.59%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '580.66'
$ws.Range('E6').Value = '  -1.41%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.669.23'
$ws.Range('E7').Value = '  +2.86%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.624'
$ws.Range('E8').Value = '  +2.35%  '

$ws.Range('E9').Value = '  +0.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.686'
$ws.Range('E10').Value = '  +2.54%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  +8.96%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '57.51'
$ws.Range('E12').Value = '  +6.80%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000299'
$ws.Range('E13').Value = '  +18.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.22'
$ws.Range('E14').Value = '  +4.47%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.248.57'
$ws.Range('E15').Value = '  +2.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.666.18'
$ws.Range('E16').Value = '  +2.50%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.127'
$ws.Range('E17').Value = '  +0.87%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.65'
$ws.Range('E18').Value = '  +4.42%  '

$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.79'
$ws.Range('E19').Value = '  +2.68%  '

$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.736.56'
$ws.Range('E20').Value = '  +3.48%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.10'
$ws.Range('E21').Value = '  +4.34%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '407.29'
$ws.Range('E22').Value = '  +3.79%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.17'
$ws.Range('E23').Value = '  +28.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.29'
$ws.Range('E24').Value = '  +0.07%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.49'
$ws.Range('E25').Value = '  +2.26%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.99'
$ws.Range('E26').Value = '  +4.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.74'
$ws.Range('E27').Value = '  +3.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.90'
$ws.Range('E28').Value = '  +9.01%  '

$ws.Range('E29').Value = '  +1.83%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.35'
$ws.Range('E30').Value = '  +23.49%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.32'
$ws.Range('E31').Value = '  +4.58%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.30'
$ws.Range('E32').Value = '  +4.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '697.46'
$ws.Range('E33').Value = '  +15.00%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.37'
$ws.Range('E34').Value = '  +3.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.118'
$ws.Range('E35').Value = '  +5.27%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '65.07'
$ws.Range('E36').Value = '  -0.25%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '43.28'
$ws.Range('E37').Value = '  +4.59%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.431'
$ws.Range('E38').Value = '  +15.74%  '

$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0811'
$ws.Range('E39').Value = '  +9.25%  '

$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('E41').Value = '  +10.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.91'
$ws.Range('E42').Value = '  +21.07%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.18'
$ws.Range('E43').Value = '  +15.15%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.234.78'
$ws.Range('E44').Value = '  +10.85%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.05'
$ws.Range('E45').Value = '  +38.04%  '

$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.997'
$ws.Range('E46').Value = '  -0.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0426'
$ws.Range('E47').Value = '  +4.75%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.02'
$ws.Range('E48').Value = '  +9.37%  '

$ws.Range('E49').Value = '  +2.60%  '

$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.65'
$ws.Range('E50').Value = '  +5.55%  '

$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.11'
$ws.Range('E51').Value = '  +2.06%  '
